$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1019.818874306688
$ws.Range("C3").Value = 993.0322683882054
$ws.Range("C4").Value = 1005.216874844783
$ws.Range("C5").Value = 1021.007972662469
$ws.Range("C6").Value = 1024.834302470246
$ws.Range("C7").Value = 1031.341892657855
$ws.Range("C8").Value = 1041.3848976627
$ws.Range("C9").Value = 1038.21761213041
$ws.Range("C10").Value = 1033.794401743227
$ws.Range("C11").Value = 1039.003312756637
$ws.Range("C12").Value = 1044.109636529123
$ws.Range("C13").Value = 1042.986982193184
$ws.Range("C14").Value = 1043.469777271211
$ws.Range("C15").Value = 1041.229210730021
$ws.Range("C16").Value = 1042.584822950721
$ws.Range("C17").Value = 1037.664461810195
$ws.Range("C18").Value = 1044.195402046921
$ws.Range("C19").Value = 1045.554744999677
$ws.Range("C20").Value = 1045.60244847629
$ws.Range("C21").Value = 1046.503106510221
$ws.Range("C22").Value = 1046.217385498926
$ws.Range("C23").Value = 1045.808824293615
$ws.Range("C24").Value = 1045.141953470607
$ws.Range("C25").Value = 1044.590254513755
$ws.Range("C26").Value = 1043.594981281037
$ws.Range("C27").Value = 1042.660623550184
$ws.Range("C28").Value = 1041.897366624411
$ws.Range("C29").Value = 1041.058522124883
$ws.Range("C30").Value = 1040.190002823183
$ws.Range("C31").Value = 1039.191022892516
$ws.Range("C32").Value = 1038.290798009989
$ws.Range("C33").Value = 1037.499605958114
$ws.Range("C34").Value = 1036.608961371366
$ws.Range("C35").Value = 1035.803462762885
$ws.Range("C36").Value = 1034.985294089741
$ws.Range("C37").Value = 1034.133429667738
$ws.Range("C38").Value = 1033.328944796685
$ws.Range("C39").Value = 1032.692957783854
$ws.Range("C40").Value = 1031.745839465081
$ws.Range("C41").Value = 1031.204667379715
$ws.Range("C42").Value = 1030.417902951775
$ws.Range("C43").Value = 1029.777211583927
$ws.Range("C44").Value = 1029.179597094809
$ws.Range("C45").Value = 1028.712653378894
$ws.Range("C46").Value = 1028.226592689419
$ws.Range("C47").Value = 1027.832577689014
$ws.Range("C48").Value = 1027.500654401603
$ws.Range("C49").Value = 1027.364752165472
$ws.Range("C50").Value = 1027.30991572586
$ws.Range("C51").Value = 1027.492102115945
$ws.Range("C52").Value = 1027.800580438054
$ws.Range("C53").Value = 1028.028716521523
$ws.Range("C54").Value = 1028.292951570253
$ws.Range("C55").Value = 1028.264552629832
$ws.Range("C56").Value = 1027.705794996225
$ws.Range("C57").Value = 1026.873158059161
$ws.Range("C58").Value = 1026.120181843115
$ws.Range("C59").Value = 1025.398385591236
$ws.Range("C60").Value = 1020.147880849143
$ws.Range("C61").Value = 1012.797518921471
$ws.Range("C62").Value = 1012.183449496384
$ws.Range("C63").Value = 1011.52381584757
$ws.Range("C64").Value = 1010.947735958223
$ws.Range("C65").Value = 1010.528340106842
